# Update "想去人数" (F column) values on the "展览" sheet (sheet1)
# and the "全部类型" sheet (sheet4) to reflect the latest scrape.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 8707
$ws1.Range("F4").Value = 38
$ws1.Range("F6").Value = 499
$ws1.Range("F7").Value = 181
$ws1.Range("F9").Value = 476
$ws1.Range("F10").Value = 84
$ws1.Range("F11").Value = 94
$ws1.Range("F13").Value = 6269
$ws1.Range("F15").Value = 327
$ws1.Range("F16").Value = 2412
$ws1.Range("F17").Value = 124
$ws1.Range("F18").Value = 213
$ws1.Range("F20").Value = 474

# --- Sheet "全部类型" (all types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 8707
$ws4.Range("F4").Value = 38
$ws4.Range("F8").Value = 499
$ws4.Range("F9").Value = 181
$ws4.Range("F11").Value = 476
$ws4.Range("F12").Value = 84
$ws4.Range("F13").Value = 94
$ws4.Range("F16").Value = 6269
$ws4.Range("F19").Value = 327
$ws4.Range("F20").Value = 2412
$ws4.Range("F21").Value = 125
$ws4.Range("F22").Value = 214
$ws4.Range("F24").Value = 474
